$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column H ("Absent") values change from 0 to 1 for rows 3 through 18
for ($row = 3; $row -le 18; $row++) {
    $ws.Cells.Item($row, 8).Value = 1
}
